$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowBF2 = New-Object 'object[,]' 1,5
$rowBF2[0,0] = 1.02
$rowBF2[0,1] = 1.064403745703216
$rowBF2[0,2] = 1.068854331059169
$rowBF2[0,3] = 1.059360194672778
$rowBF2[0,4] = 1.074484951601651
$ws.Range("B2:F2").Value = $rowBF2

$rowIN2 = New-Object 'object[,]' 1,6
$rowIN2[0,0] = 1.028616942456438
$rowIN2[0,1] = 1.069364162928255
$rowIN2[0,2] = 1.071558399391786
$rowIN2[0,3] = 1.062089911287122
$rowIN2[0,4] = 1.077174040345946
$rowIN2[0,5] = 1.070882782697397
$ws.Range("I2:N2").Value = $rowIN2

$rowBF3 = New-Object 'object[,]' 1,5
$rowBF3[0,0] = 1.02
$rowBF3[0,1] = 1.067432126304425
$rowBF3[0,2] = 1.071637106053207
$rowBF3[0,3] = 1.062016139679305
$rowBF3[0,4] = 1.077261111564872
$ws.Range("B3:F3").Value = $rowBF3

$rowIN3 = New-Object 'object[,]' 1,6
$rowIN3[0,0] = 1.028706295460768
$rowIN3[0,1] = 1.072037906608969
$rowIN3[0,2] = 1.074152358947299
$rowIN3[0,3] = 1.06455541776405
$rowIN3[0,4] = 1.07976253416825
$rowIN3[0,5] = 1.07356032340082
$ws.Range("I3:N3").Value = $rowIN3

$rowBF4 = New-Object 'object[,]' 1,5
$rowBF4[0,0] = 1.02
$rowBF4[0,1] = 1.06937848623887
$rowBF4[0,2] = 1.073425025095345
$rowBF4[0,3] = 1.063722623200556
$rowBF4[0,4] = 1.079044006668379
$ws.Range("B4:F4").Value = $rowBF4

$rowIN4 = New-Object 'object[,]' 1,6
$rowIN4[0,0] = 1.02876060958436
$rowIN4[0,1] = 1.073755061722759
$rowIN4[0,2] = 1.075817889811181
$rowIN4[0,3] = 1.066138444270211
$rowIN4[0,4] = 1.081423779036481
$rowIN4[0,5] = 1.075279917071833
$ws.Range("I4:N4").Value = $rowIN4

$rowBF5 = New-Object 'object[,]' 1,5
$rowBF5[0,0] = 1.02
$rowBF5[0,1] = 1.070193664401052
$rowBF5[0,2] = 1.074173701984362
$rowBF5[0,3] = 1.064437212360323
$rowBF5[0,4] = 1.079790392869517
$ws.Range("B5:F5").Value = $rowBF5

$rowIN5 = New-Object 'object[,]' 1,6
$rowIN5[0,0] = 1.028782609179578
$rowIN5[0,1] = 1.074473937093221
$rowIN5[0,2] = 1.076515060393026
$rowIN5[0,3] = 1.066801073486966
$rowIN5[0,4] = 1.082118967797181
$rowIN5[0,5] = 1.075999813327776
$ws.Range("I5:N5").Value = $rowIN5

$rowBF6 = New-Object 'object[,]' 1,5
$rowBF6[0,0] = 1.02
$rowBF6[0,1] = 1.07033035885954
$rowBF6[0,2] = 1.074299236641826
$rowBF6[0,3] = 1.064557032101697
$rowBF6[0,4] = 1.079915532418091
$ws.Range("B6:F6").Value = $rowBF6

$rowIN6 = New-Object 'object[,]' 1,6
$rowIN6[0,0] = 1.02878625423059
$rowIN6[0,1] = 1.074594464837127
$rowIN6[0,2] = 1.076631943626964
$rowIN6[0,3] = 1.066912165389791
$rowIN6[0,4] = 1.082235507702484
$rowIN6[0,5] = 1.076120512234904
$ws.Range("I6:N6").Value = $rowIN6

$rowBF7 = New-Object 'object[,]' 1,5
$rowBF7[0,0] = 1.02
$rowBF7[0,1] = 1.069389390639929
$rowBF7[0,2] = 1.073435040492399
$rowBF7[0,3] = 1.063732182539764
$rowBF7[0,4] = 1.079053992159789
$ws.Range("B7:F7").Value = $rowBF7

$rowIN7 = New-Object 'object[,]' 1,6
$rowIN7[0,0] = 1.028760906815136
$rowIN7[0,1] = 1.073764679117919
$rowIN7[0,2] = 1.075827217196588
$rowIN7[0,3] = 1.066147309552951
$rowIN7[0,4] = 1.081433080648688
$rowIN7[0,5] = 1.075289548124797
$ws.Range("I7:N7").Value = $rowIN7

$rowBF8 = New-Object 'object[,]' 1,5
$rowBF8[0,0] = 1.02
$rowBF8[0,1] = 1.065429989855385
$rowBF8[0,2] = 1.069797466781066
$rowBF8[0,3] = 1.060260334605696
$rowBF8[0,4] = 1.07542600533518
$ws.Range("B8:F8").Value = $rowBF8

$rowIN8 = New-Object 'object[,]' 1,6
$rowIN8[0,0] = 1.028647867919017
$rowIN8[0,1] = 1.070270494411261
$rowIN8[0,2] = 1.072437765768965
$rowIN8[0,3] = 1.062925737286188
$rowIN8[0,4] = 1.078051714336241
$rowIN8[0,5] = 1.071790401275067
$ws.Range("I8:N8").Value = $rowIN8

$rowBF9 = New-Object 'object[,]' 1,5
$rowBF9[0,0] = 1.02
$rowBF9[0,1] = 1.058347777830567
$rowBF9[0,2] = 1.063286424764229
$rowBF9[0,3] = 1.054046324453896
$rowBF9[0,4] = 1.068926189156293
$ws.Range("B9:F9").Value = $rowBF9

$rowIN9 = New-Object 'object[,]' 1,6
$rowIN9[0,0] = 1.02842163407217
$rowIN9[0,1] = 1.064010546241992
$rowIN9[0,2] = 1.066362486643547
$rowIN9[0,3] = 1.057151135187785
$rowIN9[0,4] = 1.071984970469591
$rowIN9[0,5] = 1.065521563261371
$ws.Range("I9:N9").Value = $rowIN9

$rowBF10 = New-Object 'object[,]' 1,5
$rowBF10[0,0] = 1.02
$rowBF10[0,1] = 1.053549684627253
$rowBF10[0,2] = 1.05887238238032
$rowBF10[0,3] = 1.049833897324116
$rowBF10[0,4] = 1.064515858065685
$ws.Range("B10:F10").Value = $rowBF10

$rowIN10 = New-Object 'object[,]' 1,6
$rowIN10[0,0] = 1.028252318237376
$rowIN10[0,1] = 1.059762922972372
$rowIN10[0,2] = 1.062238202889539
$rowIN10[0,3] = 1.053230786391808
$rowIN10[0,4] = 1.067862573882963
$rowIN10[0,5] = 1.061267907879499
$ws.Range("I10:N10").Value = $rowIN10

$rowBF11 = New-Object 'object[,]' 1,5
$rowBF11[0,0] = 1.02
$rowBF11[0,1] = 1.051452545755849
$rowBF11[0,2] = 1.056942439544046
$rowBF11[0,3] = 1.047992161441327
$rowBF11[0,4] = 1.062586635767228
$ws.Range("B11:F11").Value = $rowBF11

$rowIN11 = New-Object 'object[,]' 1,6
$rowIN11[0,0] = 1.028174544183082
$rowIN11[0,1] = 1.057904828556596
$rowIN11[0,2] = 1.060433603597233
$rowIN11[0,3] = 1.051515370051791
$rowIN11[0,4] = 1.066057897873709
$rowIN11[0,5] = 1.059407174756527
$ws.Range("I11:N11").Value = $rowIN11

$rowBF12 = New-Object 'object[,]' 1,5
$rowBF12[0,0] = 1.02
$rowBF12[0,1] = 1.050670522904449
$rowBF12[0,2] = 1.056222666572231
$rowBF12[0,3] = 1.047305293258374
$rowBF12[0,4] = 1.061866997945824
$ws.Range("B12:F12").Value = $rowBF12

$rowIN12 = New-Object 'object[,]' 1,6
$rowIN12[0,0] = 1.028144979039361
$rowIN12[0,1] = 1.057211713051172
$rowIN12[0,2] = 1.059760374797199
$rowIN12[0,3] = 1.050875404097188
$rowIN12[0,4] = 1.065384506959405
$rowIN12[0,5] = 1.058713074947583
$ws.Range("I12:N12").Value = $rowIN12

$rowBF13 = New-Object 'object[,]' 1,5
$rowBF13[0,0] = 1.02
$rowBF13[0,1] = 1.050838409500777
$rowBF13[0,2] = 1.056377193583822
$rowBF13[0,3] = 1.047452755633348
$rowBF13[0,4] = 1.062021501951339
$ws.Range("B13:F13").Value = $rowBF13

$rowIN13 = New-Object 'object[,]' 1,6
$rowIN13[0,0] = 1.028151351577729
$rowIN13[0,1] = 1.057360523287193
$rowIN13[0,2] = 1.059904918484111
$rowIN13[0,3] = 1.051012806548801
$rowIN13[0,4] = 1.065529091455907
$rowIN13[0,5] = 1.058862096511208
$ws.Range("I13:N13").Value = $rowIN13

$rowBF14 = New-Object 'object[,]' 1,5
$rowBF14[0,0] = 1.02
$rowBF14[0,1] = 1.051387966457541
$rowBF14[0,2] = 1.056883002802341
$rowBF14[0,3] = 1.047935441720053
$rowBF14[0,4] = 1.062527212895788
$ws.Range("B14:F14").Value = $rowBF14

$rowIN14 = New-Object 'object[,]' 1,6
$rowIN14[0,0] = 1.02817211415558
$rowIN14[0,1] = 1.057847595936558
$rowIN14[0,2] = 1.06037801446979
$rowIN14[0,3] = 1.051462527634266
$rowIN14[0,4] = 1.066002298067677
$rowIN14[0,5] = 1.059349860859603
$ws.Range("I14:N14").Value = $rowIN14

$rowBF15 = New-Object 'object[,]' 1,5
$rowBF15[0,0] = 1.02
$rowBF15[0,1] = 1.05172615862937
$rowBF15[0,2] = 1.057194260197556
$rowBF15[0,3] = 1.04823247098932
$rowBF15[0,4] = 1.062838392205933
$ws.Range("B15:F15").Value = $rowBF15

$rowIN15 = New-Object 'object[,]' 1,6
$rowIN15[0,0] = 1.028184816838807
$rowIN15[0,1] = 1.058147305156268
$rowIN15[0,2] = 1.060669114447236
$rowIN15[0,3] = 1.051739243764195
$rowIN15[0,4] = 1.066293448520204
$rowIN15[0,5] = 1.059649995700783
$ws.Range("I15:N15").Value = $rowIN15

$rowBF16 = New-Object 'object[,]' 1,5
$rowBF16[0,0] = 1.02
$rowBF16[0,1] = 1.053688443340699
$rowBF16[0,2] = 1.059000064700035
$rowBF16[0,3] = 1.049955745125523
$rowBF16[0,4] = 1.064643473947556
$ws.Range("B16:F16").Value = $rowBF16

$rowIN16 = New-Object 'object[,]' 1,6
$rowIN16[0,0] = 1.028257385349631
$rowIN16[0,1] = 1.059885832618036
$rowIN16[0,2] = 1.06235756435186
$rowIN16[0,3] = 1.053344247945791
$rowIN16[0,4] = 1.067981921661329
$rowIN16[0,5] = 1.061390992070958
$ws.Range("I16:N16").Value = $rowIN16

$rowBF17 = New-Object 'object[,]' 1,5
$rowBF17[0,0] = 1.02
$rowBF17[0,1] = 1.054914021315716
$rowBF17[0,2] = 1.06012773504143
$rowBF17[0,3] = 1.051031892207867
$rowBF17[0,4] = 1.065770454070924
$ws.Range("B17:F17").Value = $rowBF17

$rowIN17 = New-Object 'object[,]' 1,6
$rowIN17[0,0] = 1.028301707324655
$rowIN17[0,1] = 1.06097124619261
$rowIN17[0,2] = 1.063411591307194
$rowIN17[0,3] = 1.054346169708111
$rowIN17[0,4] = 1.069035724351104
$rowIN17[0,5] = 1.062477947057307
$ws.Range("I17:N17").Value = $rowIN17

$rowBF18 = New-Object 'object[,]' 1,5
$rowBF18[0,0] = 1.02
$rowBF18[0,1] = 1.055627003683143
$rowBF18[0,2] = 1.060783695593508
$rowBF18[0,3] = 1.051657887230524
$rowBF18[0,4] = 1.066425926451106
$ws.Range("B18:F18").Value = $rowBF18

$rowIN18 = New-Object 'object[,]' 1,6
$rowIN18[0,0] = 1.028327129784685
$rowIN18[0,1] = 1.061602538505062
$rowIN18[0,2] = 1.064024584650605
$rowIN18[0,3] = 1.054928855358737
$rowIN18[0,4] = 1.069648500592482
$rowIN18[0,5] = 1.063110135877253
$ws.Range("I18:N18").Value = $rowIN18

$rowBF19 = New-Object 'object[,]' 1,5
$rowBF19[0,0] = 1.02
$rowBF19[0,1] = 1.055869797463656
$rowBF19[0,2] = 1.061007060685006
$rowBF19[0,3] = 1.051871049488389
$rowBF19[0,4] = 1.066649110551965
$ws.Range("B19:F19").Value = $rowBF19

$rowIN19 = New-Object 'object[,]' 1,6
$rowIN19[0,0] = 1.028335725477443
$rowIN19[0,1] = 1.061817488824176
$rowIN19[0,2] = 1.064233296815799
$rowIN19[0,3] = 1.055127247564125
$rowIN19[0,4] = 1.069857124074533
$rowIN19[0,5] = 1.063325391450475
$ws.Range("I19:N19").Value = $rowIN19

$rowBF20 = New-Object 'object[,]' 1,5
$rowBF20[0,0] = 1.02
$rowBF20[0,1] = 1.054782723334973
$rowBF20[0,2] = 1.060006932646272
$rowBF20[0,3] = 1.050916608726803
$rowBF20[0,4] = 1.065649734576553
$ws.Range("B20:F20").Value = $rowBF20

$rowIN20 = New-Object 'object[,]' 1,6
$rowIN20[0,0] = 1.028296996498232
$rowIN20[0,1] = 1.060854979665842
$rowIN20[0,2] = 1.063298691386386
$rowIN20[0,3] = 1.054238851421313
$rowIN20[0,4] = 1.068922857412594
$rowIN20[0,5] = 1.062361515418735
$ws.Range("I20:N20").Value = $rowIN20

$rowBF21 = New-Object 'object[,]' 1,5
$rowBF21[0,0] = 1.02
$rowBF21[0,1] = 1.051226220844427
$rowBF21[0,2] = 1.056734135694814
$rowBF21[0,3] = 1.047793379873758
$rowBF21[0,4] = 1.062378378379244
$ws.Range("B21:F21").Value = $rowBF21

$rowIN21 = New-Object 'object[,]' 1,6
$rowIN21[0,0] = 1.028166018818732
$rowIN21[0,1] = 1.057704247102953
$rowIN21[0,2] = 1.060238780942868
$rowIN21[0,3] = 1.051330173633192
$rowIN21[0,4] = 1.065863035647157
$rowIN21[0,5] = 1.059206308454213
$ws.Range("I21:N21").Value = $rowIN21

$rowBF22 = New-Object 'object[,]' 1,5
$rowBF22[0,0] = 1.02
$rowBF22[0,1] = 1.04897238736131
$rowBF22[0,2] = 1.054659528384192
$rowBF22[0,3] = 1.045813628986947
$rowBF22[0,4] = 1.060303910738075
$ws.Range("B22:F22").Value = $rowBF22

$rowIN22 = New-Object 'object[,]' 1,6
$rowIN22[0,0] = 1.028079751163292
$rowIN22[0,1] = 1.055706212423497
$rowIN22[0,2] = 1.058297945196023
$rowIN22[0,3] = 1.04948521490176
$rowIN22[0,4] = 1.063921482925037
$rowIN22[0,5] = 1.057205436336334
$ws.Range("I22:N22").Value = $rowIN22

$rowBF23 = New-Object 'object[,]' 1,5
$rowBF23[0,0] = 1.02
$rowBF23[0,1] = 1.050168906849766
$rowBF23[0,2] = 1.055760952536984
$rowBF23[0,3] = 1.046864688921008
$rowBF23[0,4] = 1.061405333240324
$ws.Range("B23:F23").Value = $rowBF23

$rowIN23 = New-Object 'object[,]' 1,6
$rowIN23[0,0] = 1.02812585676196
$rowIN23[0,1] = 1.056767059920024
$rowIN23[0,2] = 1.05932846039377
$rowIN23[0,3] = 1.050464827376597
$rowIN23[0,4] = 1.064952451222203
$rowIN23[0,5] = 1.058267790357986
$ws.Range("I23:N23").Value = $rowIN23

$rowBF24 = New-Object 'object[,]' 1,5
$rowBF24[0,0] = 1.02
$rowBF24[0,1] = 1.054842057032564
$rowBF24[0,2] = 1.060061523575272
$rowBF24[0,3] = 1.050968705627154
$rowBF24[0,4] = 1.065704288310773
$ws.Range("B24:F24").Value = $rowBF24

$rowIN24 = New-Object 'object[,]' 1,6
$rowIN24[0,0] = 1.028299126445682
$rowIN24[0,1] = 1.060907521094832
$rowIN24[0,2] = 1.063349711565542
$rowIN24[0,3] = 1.054287349239813
$rowIN24[0,4] = 1.068973862955187
$rowIN24[0,5] = 1.06241413146258
$ws.Range("I24:N24").Value = $rowIN24

$rowBF25 = New-Object 'object[,]' 1,5
$rowBF25[0,0] = 1.02
$rowBF25[0,1] = 1.06019176451774
$rowBF25[0,2] = 1.064982213898239
$rowBF25[0,3] = 1.055664707598433
$rowBF25[0,4] = 1.070619743892874
$ws.Range("B25:F25").Value = $rowBF25

$rowIN25 = New-Object 'object[,]' 1,6
$rowIN25[0,0] = 1.028483357108061
$rowIN25[0,1] = 1.065641599946277
$rowIN25[0,2] = 1.067945770163264
$rowIN25[0,3] = 1.058656089067854
$rowIN25[0,4] = 1.073566720477463
$rowIN25[0,5] = 1.070882782697397
$ws.Range("I25:N25").Value = $rowIN25

Write-Host "Applied vm_pu updates for 380 kV case"